$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPIs")

# Fix template weights: swap the Weight values for rows 5 and 6 (H5 <-> H6)
$ws.Range("H5").Value = 0.15
$ws.Range("H6").Value = 0.35

# Update selection/view state to match the saved workbook
$ws.Range("H7").Select()
$excel.ActiveWindow.ScrollColumn = 4
